# Add "2022-Q3" quarter data:
#  1. Insert a new worksheet named "2022-Q3" right after "总计" (shifts
#     2022-Q2 / 2021-Q3 / 2021-Q2 / 2021-Q1 / 2020-Q4 one slot to the right).
#  2. Fill it with the fund-holdings table for that quarter.
#  3. Prepend a matching summary row to the "总计" sheet and shift its
#     existing rows down by one.

$wb = $excel.ActiveWorkbook
$zongji = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create + position the new sheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Helper: write a value that looks numeric but must stay TEXT (matches the
# inlineStr cells used for fund codes / percentages in sibling sheets).
# Leading "'" forces Excel to keep it as text; resetting Style to "Normal"
# afterwards drops the quotePrefix formatting bit it would otherwise pick up.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2. Header row (B1:H1) - text, styled like every other quarter sheet
# ---------------------------------------------------------------------
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

$zongji.Cells.Item(1,2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122, $null, $false, $false)

# ---------------------------------------------------------------------
# 3. Data rows 2-16
# ---------------------------------------------------------------------
$rows = @(
    @(0,  "015182", "汇添富逆向投资混合D",              "21.25", "92.44", "4.26", "0.9052", 5),
    @(1,  "470098", "汇添富逆向投资混合A",              "20.96", "92.44", "4.26", "0.8929", 5),
    @(2,  "260116", "景顺长城核心竞争力混合A",          "16.26", "72.95", "4.06", "0.6602", 4),
    @(3,  "009190", "景顺长城核心优选一年持有期混合",    "11.52", "86.01", "4.88", "0.5622", 5),
    @(4,  "011346", "淳厚鑫淳一年持有期混合",            "4.81",  "69.72", "4.49", "0.2160", 1),
    @(5,  "012454", "淳厚鑫悦混合A",                    "2.06",  "75.61", "5.39", "0.1110", 1),
    @(6,  "001541", "汇添富民营新动力股票",              "2.21",  "91.50", "3.58", "0.0791", 5),
    @(7,  "014509", "汇添富先进制造混合C",               "0.85",  "87.50", "4.70", "0.0400", 8),
    @(8,  "014508", "汇添富先进制造混合A",               "0.83",  "87.50", "4.70", "0.0390", 8),
    @(9,  "012455", "淳厚鑫悦混合C",                     "0.68",  "75.61", "5.39", "0.0367", 1),
    @(10, "960008", "景顺长城核心竞争力混合H",           "0.34",  "72.95", "4.06", "0.0138", 4),
    @(11, "015181", "汇添富逆向投资混合C",               "0.24",  "92.44", "4.26", "0.0102", 5),
    @(12, "013368", "汇添富多元价值发现混合C",           "0.28",  "36.00", "1.73", "0.0048", 3),
    @(13, "013367", "汇添富多元价值发现混合A",           "0.25",  "36.00", "1.73", "0.0043", 3),
    @(14, "015731", "景顺长城核心竞争力混合C",           "0.04",  "72.95", "4.06", "0.0016", 4)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q3.Cells.Item($r, 2) $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3.Cells.Item($r, 6) $row[5]
    Set-TextValue $q3.Cells.Item($r, 7) $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Column-A style (s="2") matches every sibling quarter sheet
$zongji.Cells.Item(2,1).Copy()
$q3.Range("A2:A16").PasteSpecial(-4122, $null, $false, $false)

Write-Host "2022-Q3 sheet populated"
